{"js": "// Obituary.docx update:\n//  1. Merge the split runs in three \"Content\" cells (\"AKC Canine  Health\n//     Foundation\", \"LinkedIn  Profile\", \"Photography Portfolio\") into a\n//     single run each (removing the stray extra spaces/run-splits).\n//  2. Append a new table row \"GitHub Repositories\" ->\n//     https://github.com/RalphHightower (as a real hyperlink, matching\n//     the style of the other rows).\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// --- 1. Clean up the three split-run cells in column 0 -----------------\nconst fixes = [\n  { row: 1, text: \"AKC Canine Health Foundation\" },\n  { row: 2, text: \"LinkedIn Profile\" },\n  { row: 3, text: \"Photography Portfolio\" },\n];\n\nfor (const fix of fixes) {\n  const cell = table.getCell(fix.row, 0);\n  const para = cell.body.paragraphs.getFirst();\n  // Replace only the paragraph's content (keeps the Compact/left pPr\n  // formatting intact) and collapses every run down to one.\n  para.getRange(\"Content\").insertText(fix.text, \"Replace\");\n}\nawait context.sync();\n\n// --- 2. Add the new \"GitHub Repositories\" row ---------------------------\ntable.addRows(\"End\", 1, [[\"GitHub Repositories\", \"https://github.com/RalphHightower\"]]);\nawait context.sync();\n\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst newRowIndex = table.rowCount - 1;\nconst linkCell = table.getCell(newRowIndex, 1);\nconst linkPara = linkCell.body.paragraphs.getFirst();\nconst linkRange = linkPara.getRange();\nlinkRange.hyperlink = \"https://github.com/RalphHightower\";\nawait context.sync();\n", "ps1": "# Obituary.docx update:\n#  1. Merge the split runs in three \"Content\" cells (\"AKC Canine  Health\n#     Foundation\", \"LinkedIn  Profile\", \"Photography Portfolio\") into a\n#     single run each (removing the stray extra spaces/run-splits).\n#  2. Append a new table row \"GitHub Repositories\" ->\n#     https://github.com/RalphHightower (as a real hyperlink, matching\n#     the style of the other rows).\n\n$d = $word.ActiveDocument\n\n# --- 1. Clean up the three split-run cells in column 0 --------------------\n$fixes = @(\n    @{ Find = \"AKC Canine  Health Foundation\"; Replace = \"AKC Canine Health Foundation\" },\n    @{ Find = \"LinkedIn  Profile\";             Replace = \"LinkedIn Profile\" },\n    @{ Find = \"Photography Portfolio\";         Replace = \"Photography Portfolio\" }\n)\n\nforeach ($fix in $fixes) {\n    $rng = $d.Content\n    $rng.Find.Execute($fix.Find, $false, $false, $false, $false, $false, $true, 1, $false, $fix.Replace, 2)\n}\n\n# --- 2. Add the new \"GitHub Repositories\" row ------------------------------\n$table = $d.Tables(1)\n$table.Rows.Add() | Out-Null\n$rowCount = $table.Rows.Count\n\n$labelCell = $table.Cell($rowCount, 1)\n$labelCell.Range.Text = \"GitHub Repositories\"\n\n$linkText = \"https://github.com/RalphHightower\"\n$linkCell = $table.Cell($rowCount, 2)\n$linkCell.Range.Text = $linkText\n\n# Locate the freshly inserted URL text as its own sub-range so the\n# hyperlink wraps exactly that text (avoids leaving a stray empty run\n# behind in the cell).\n$linkRange = $linkCell.Range.Duplicate\n$linkRange.Find.Execute($linkText) | Out-Null\n$d.Hyperlinks.Add($linkRange, $linkText) | Out-Null\n"}
